# Handback status report: refresh the generated timestamps.
#
# "Latest HO Xliff Generate Date" (Overview!G2) and "Correspond Handoff
# Datetime" (de-de!H2) share the same underlying value
# ("2016-09-07 17:35:56" -> "2016-09-07 17:37:14"), so both cells are set
# to the new value together.
#
# The zh-cn sheet's "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) each get their own refreshed
# timestamp, as does de-de's "Correspond Handback DateTime" (K2).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview) / Correspond Handoff Datetime (de-de)
$wsOverview.Range("G2").Value = "2016-09-07 17:37:14"
$wsDeDe.Range("H2").Value = "2016-09-07 17:37:14"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-07 17:37:01"
$wsZhCn.Range("K2").Value = "2016-09-07 17:37:32"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-07 17:37:41"
